# Update "Hoja1" sheet A1 text with new conversion rates
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 11.51 = 47196.85 pesos`n✅ 47196.85 pesos = 11.49 = 976.33 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Update "tasas" sheet N10, O10, N12, O12 values
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 86.87
$ws2.Range("O10").Value = 4099.99
$ws2.Range("N12").Value = 4109
$ws2.Range("O12").Value = 85
